$wb = $excel.ActiveWorkbook

# Sheet "展览": update 想去人数 (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13939
$ws1.Range("F4").Value = 672
$ws1.Range("F6").Value = 515
$ws1.Range("F7").Value = 1446

# Sheet "全部类型": update 想去人数 (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13939
$ws4.Range("F4").Value = 672
$ws4.Range("F8").Value = 515
$ws4.Range("F9").Value = 1446
